$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns retain their original text
# representation (e.g. "128.00", "0.0000130") instead of being
# auto-converted to numbers by Excel when the value looks numeric.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '53.981.11'
$ws.Range("E2").Value = '  -2.67%  '
$ws.Range("D3").Value = '2.287.91'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '495.48'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").Value = '128.00'
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  -1.25%  '
$ws.Range("D9").Value = '2.287.43'
$ws.Range("E9").Value = '  -1.78%  '
$ws.Range("D10").Value = '0.0940'
$ws.Range("E10").Value = '  -3.50%  '
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("E13").Value = '  -3.00%  '
$ws.Range("D14").Value = '2.692.53'
$ws.Range("E14").Value = '  -2.00%  '
$ws.Range("D15").Value = '21.56'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").Value = '54.022.11'
$ws.Range("E16").Value = '  -2.60%  '
$ws.Range("D17").Value = '0.0000130'
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").Value = '2.297.87'
$ws.Range("E18").Value = '  -4.13%  '
$ws.Range("D19").Value = '9.94'
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("E20").Value = '  +1.48%  '
$ws.Range("D21").Value = '299.46'
$ws.Range("E21").Value = '  -2.66%  '
$ws.Range("D22").Value = '6.28'
$ws.Range("E22").Value = '  +1.46%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '63.54'
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = '0.375'
$ws.Range("E26").Value = '  +1.47%  '
$ws.Range("D27").Value = '2.385.93'
$ws.Range("E27").Value = '  -3.51%  '
$ws.Range("E28").Value = '  +1.64%  '
$ws.Range("D29").Value = '7.16'
$ws.Range("E29").Value = '  +0.82%  '
$ws.Range("D30").Value = '162.99'
$ws.Range("E30").Value = '  -5.65%  '
$ws.Range("D31").Value = '1.61'
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("D32").Value = '0.0₃0685'
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  -0.37%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("D37").Value = '17.51'
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("E38").Value = '  +0.99%  '
$ws.Range("D39").Value = '0.863'
$ws.Range("E39").Value = '  +5.48%  '
$ws.Range("E40").Value = '  +0.89%  '
$ws.Range("D41").Value = '35.40'
$ws.Range("E41").Value = '  -1.86%  '
$ws.Range("D42").Value = '0.376'
$ws.Range("E42").Value = '  +1.65%  '
$ws.Range("E43").Value = '  +1.80%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = '4.92'
$ws.Range("E45").Value = '  +3.51%  '
$ws.Range("D46").Value = '126.83'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").Value = '0.0891'
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").Value = '0.549'
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("D49").Value = '239.29'
$ws.Range("E49").Value = '  +1.80%  '
$ws.Range("E50").Value = '  +1.05%  '
$ws.Range("E51").Value = '  -0.25%  '
